# Updated cryptos list - applying price/volume/coin changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.206.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.242.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.33%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.73%  '
$ws.Range('E6').Value = '  -5.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '69.20'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.32%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.563'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.83%  '
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.03'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '35.30'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.58%  '
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.581.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.96%  '
$ws.Range('E17').Value = '  -5.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.248.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.129.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('E21').Value = '  -6.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.86%  '
$ws.Range('E23').Value = '  -6.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.52'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.94%  '
$ws.Range('E32').Value = '  -5.55%  '
$ws.Range('E33').Value = '  -6.37%  '
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.31'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.71'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.94%  '
$ws.Range('E38').Value = '  -5.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '20.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.57%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0266'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('B45').Value = 'BitTorrent-New'
$ws.Range('C45').Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₃0165'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +30.14%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.102'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.82%  '
$ws.Range('B47').Value = 'BinanceUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.188'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.72%  '
$ws.Range('B49').Value = 'SynthetixNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.29%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.31%  '
